# Generate Report for Handback
# Removes the row-3 entries (the "ad6ac1e5-0151-40da-8053-6be78554db72" file)
# from the Overview / zh-cn / de-de sheets, and refreshes the handoff/handback
# timestamps recorded for the remaining (1117757f-...) row.

$wb = $excel.ActiveWorkbook

function Remove-Row3Hyperlinks($ws) {
    # Deleting hyperlinks while enumerating the live collection can skip
    # entries (the collection re-indexes after each delete), so restart the
    # scan after every removal instead of trying to delete mid-iteration.
    $more = $true
    while ($more) {
        $more = $false
        foreach ($hl in $ws.Hyperlinks) {
            if ($hl.Range.Row -eq 3) {
                $hl.Delete()
                $more = $true
                break
            }
        }
    }
}

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
Remove-Row3Hyperlinks $wsOverview
$wsOverview.Rows(3).Delete()

# --- zh-cn sheet --------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
Remove-Row3Hyperlinks $wsZh
$wsZh.Rows(3).Delete()
$wsZh.Range("E2").Value = "2016-03-13 04:48:37"
$wsZh.Range("H2").Value = "2016-03-13 04:48:54"

# --- de-de sheet --------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
Remove-Row3Hyperlinks $wsDe
$wsDe.Rows(3).Delete()
$wsDe.Range("E2").Value = "2016-03-13 04:48:41"
$wsDe.Range("H2").Value = "2016-03-13 04:49:00"
